$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 223
$ws.Range("I2").Value = 605
$ws.Range("J2").Value = 2339
$ws.Range("K2").Value = 11
$ws.Range("L2").Value = 693
$ws.Range("M2").Value = 57
$ws.Range("N2").Value = 446
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 8
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 29
$ws.Range("S2").Value = 266
$ws.Range("T2").Value = 438
$ws.Range("U2").Value = 29
$ws.Range("V2").Value = 3720
$ws.Range("X2").Value = 3732
$ws.Range("Y2").Value = 12
$ws.Range("AA2").Value = 22
